$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-22 (demscalar_ippu_*): set J:AS (columns 10-45) to 1
$onesRows = @(4,5,6,7,8,9,10,12,13,14,15,21,22)
foreach ($r in $onesRows) {
    for ($c = 10; $c -le 45; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# Row 96
$ws.Cells.Item(96, 10).Value = 6270808.16669041
$ws.Cells.Item(96, 11).Value = 6780525.6716
$ws.Cells.Item(96, 12).Value = 7061019.63121261
$ws.Cells.Item(96, 13).Value = 5400697.35154045
$ws.Cells.Item(96, 14).Value = 6099237.2685588
$ws.Cells.Item(96, 15).Value = 6336705.94700221
$ws.Cells.Item(96, 16).Value = 6336706
$ws.Cells.Item(96, 17).Value = 6832960.64336848
$ws.Cells.Item(96, 18).Value = 7368079.11773445
$ws.Cells.Item(96, 19).Value = 7945105.016503
$ws.Cells.Item(96, 20).Value = 8567320.28994157
$ws.Cells.Item(96, 21).Value = 9238263.911928849
$ws.Cells.Item(96, 22).Value = 9961752.008576879
$ws.Cells.Item(96, 23).Value = 10741899.563212
$ws.Cells.Item(96, 24).Value = 11583143.8211658
$ws.Cells.Item(96, 25).Value = 12490269.5274963
$ws.Cells.Item(96, 26).Value = 13468436.1411823
$ws.Cells.Item(96, 27).Value = 14523207.1805794
$ws.Cells.Item(96, 28).Value = 15660581.8670435
$ws.Cells.Item(96, 29).Value = 16887029.2467029
$ws.Cells.Item(96, 30).Value = 18209524.9844531
$ws.Cells.Item(96, 31).Value = 19635591.0394459
$ws.Cells.Item(96, 32).Value = 21173338.4477381
$ws.Cells.Item(96, 33).Value = 22831513.4554317
$ws.Cells.Item(96, 34).Value = 24619547.2646992
$ws.Cells.Item(96, 35).Value = 26547609.6756327
$ws.Cells.Item(96, 36).Value = 28626666.9290174
$ws.Cells.Item(96, 37).Value = 30868544.0790206
$ws.Cells.Item(96, 38).Value = 33285992.250553
$ws.Cells.Item(96, 39).Value = 35892761.163844
$ws.Cells.Item(96, 40).Value = 38703677.3387262
$ws.Cells.Item(96, 41).Value = 41734728.4234346
$ws.Cells.Item(96, 42).Value = 45003154.1275547
$ws.Cells.Item(96, 43).Value = 48527544.2763207
$ws.Cells.Item(96, 44).Value = 52327944.5439667
$ws.Cells.Item(96, 45).Value = 56425970.4675099

# Row 97
$ws.Cells.Item(97, 10).Value = 63860.6356247595
$ws.Cells.Item(97, 11).Value = 62774.0764733532
$ws.Cells.Item(97, 12).Value = 67726.4549433253
$ws.Cells.Item(97, 13).Value = 60733.722558311
$ws.Cells.Item(97, 14).Value = 57053.8711792603
$ws.Cells.Item(97, 15).Value = 64326.1192831146
$ws.Cells.Item(97, 16).Value = 64326.12
$ws.Cells.Item(97, 17).Value = 66740.5104137818
$ws.Cells.Item(97, 18).Value = 69245.521574939
$ws.Cells.Item(97, 19).Value = 71844.55480573
$ws.Cells.Item(97, 20).Value = 74541.13909225891
$ws.Cells.Item(97, 21).Value = 77338.9358761583
$ws.Cells.Item(97, 22).Value = 80241.744026121
$ws.Cells.Item(97, 23).Value = 83253.5049960318
$ws.Cells.Item(97, 24).Value = 86378.3081767017
$ws.Cells.Item(97, 25).Value = 89620.3964484724
$ws.Cells.Item(97, 26).Value = 92984.1719422301
$ws.Cells.Item(97, 27).Value = 96474.2020166502
$ws.Cells.Item(97, 28).Value = 100095.22545979
$ws.Cells.Item(97, 29).Value = 103852.158923451
$ws.Cells.Item(97, 30).Value = 107750.10359904
$ws.Cells.Item(97, 31).Value = 111794.352144009
$ws.Cells.Item(97, 32).Value = 115990.395868261
$ws.Cells.Item(97, 33).Value = 120343.932190288
$ws.Cells.Item(97, 34).Value = 124860.872373172
$ws.Cells.Item(97, 35).Value = 129547.349550938
$ws.Cells.Item(97, 36).Value = 134409.727056166
$ws.Cells.Item(97, 37).Value = 139454.607060175
$ws.Cells.Item(97, 38).Value = 144688.839537491
$ws.Cells.Item(97, 39).Value = 150119.531566803
$ws.Cells.Item(97, 40).Value = 155754.056980994
$ws.Cells.Item(97, 41).Value = 161600.066379393
$ws.Cells.Item(97, 42).Value = 167665.49751581
$ws.Cells.Item(97, 43).Value = 173958.586076477
$ws.Cells.Item(97, 44).Value = 180487.876862523
$ws.Cells.Item(97, 45).Value = 187262.235392164

# Row 98
$ws.Cells.Item(98, 10).Value = 62564.0860174214
$ws.Cells.Item(98, 11).Value = 70554.0525673035
$ws.Cells.Item(98, 12).Value = 77457.1410833503
$ws.Cells.Item(98, 13).Value = 93226.9265936227
$ws.Cells.Item(98, 14).Value = 91221.56263824389
$ws.Cells.Item(98, 15).Value = 90719.24332959981
$ws.Cells.Item(98, 16).Value = 90719.24000000001
$ws.Cells.Item(98, 17).Value = 97156.73433280201
$ws.Cells.Item(98, 18).Value = 104051.037312643
$ws.Cells.Item(98, 19).Value = 111434.564368451
$ws.Cells.Item(98, 20).Value = 119342.031148377
$ws.Cells.Item(98, 21).Value = 127810.616744804
$ws.Cells.Item(98, 22).Value = 136880.138501893
$ws.Cells.Item(98, 23).Value = 146593.239227595
$ws.Cells.Item(98, 24).Value = 156995.587690332
$ws.Cells.Item(98, 25).Value = 168136.093343063
$ws.Cells.Item(98, 26).Value = 180067.136284289
$ws.Cells.Item(98, 27).Value = 192844.813537251
$ws.Cells.Item(98, 28).Value = 206529.202805242
$ws.Cells.Item(98, 29).Value = 221184.644943171
$ws.Cells.Item(98, 30).Value = 236880.04647348
$ws.Cells.Item(98, 31).Value = 253689.203568787
$ws.Cells.Item(98, 32).Value = 271691.149024537
$ws.Cells.Item(98, 33).Value = 290970.523853051
$ws.Cells.Item(98, 34).Value = 311617.975246123
$ws.Cells.Item(98, 35).Value = 333730.582777295
$ws.Cells.Item(98, 36).Value = 357412.314847708
$ws.Cells.Item(98, 37).Value = 382774.517521647
$ws.Cells.Item(98, 38).Value = 409936.438050155
$ws.Cells.Item(98, 39).Value = 439025.785544213
$ws.Cells.Item(98, 40).Value = 470179.33143365
$ws.Cells.Item(98, 41).Value = 503543.552534981
$ws.Cells.Item(98, 42).Value = 539275.319751758
$ws.Cells.Item(98, 43).Value = 577542.635645519
$ws.Cells.Item(98, 44).Value = 618525.424345243
$ws.Cells.Item(98, 45).Value = 662416.377509273

# Row 99
$ws.Cells.Item(99, 10).Value = 4850.22797342753
$ws.Cells.Item(99, 11).Value = 5439.5646939778
$ws.Cells.Item(99, 12).Value = 5444.62232709263
$ws.Cells.Item(99, 13).Value = 6234.19263835195
$ws.Cells.Item(99, 14).Value = 5893.27057727725
$ws.Cells.Item(99, 15).Value = 5821.27029202642
$ws.Cells.Item(99, 16).Value = 5821.27
$ws.Cells.Item(99, 17).Value = 6025.39215786481
$ws.Cells.Item(99, 18).Value = 6236.67183553739
$ws.Cells.Item(99, 19).Value = 6455.36001062024
$ws.Cells.Item(99, 20).Value = 6681.71646121641
$ws.Cells.Item(99, 21).Value = 6916.01007451803
$ws.Cells.Item(99, 22).Value = 7158.51916621544
$ws.Cells.Item(99, 23).Value = 7409.53181110641
$ws.Cells.Item(99, 24).Value = 7669.34618529812
$ws.Cells.Item(99, 25).Value = 7938.27092040838
$ws.Cells.Item(99, 26).Value = 8216.6254701869
$ws.Cells.Item(99, 27).Value = 8504.740489992109
$ws.Cells.Item(99, 28).Value = 8802.95822957426
$ws.Cells.Item(99, 29).Value = 9111.63293963143
$ws.Cells.Item(99, 30).Value = 9431.13129262134
$ws.Cells.Item(99, 31).Value = 9761.832818329
$ws.Cells.Item(99, 32).Value = 10104.1303547073
$ws.Cells.Item(99, 33).Value = 10458.4305145264
$ws.Cells.Item(99, 34).Value = 10825.1541683862
$ws.Cells.Item(99, 35).Value = 11204.7369446652
$ws.Cells.Item(99, 36).Value = 11597.6297470009
$ws.Cells.Item(99, 37).Value = 12004.2992899142
$ws.Cells.Item(99, 38).Value = 12425.2286532168
$ws.Cells.Item(99, 39).Value = 12860.9178558578
$ws.Cells.Item(99, 40).Value = 13311.8844498931
$ws.Cells.Item(99, 41).Value = 13778.6641352812
$ws.Cells.Item(99, 42).Value = 14261.8113962377
$ws.Cells.Item(99, 43).Value = 14761.9001599028
$ws.Cells.Item(99, 44).Value = 15279.5244781056
$ws.Cells.Item(99, 45).Value = 15815.2992330334

# Row 100
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 1508.71573858184
$ws.Cells.Item(100, 12).Value = 1728.28850397102
$ws.Cells.Item(100, 13).Value = 2028.94716121149
$ws.Cells.Item(100, 14).Value = 0
$ws.Cells.Item(100, 15).Value = 560.0833316468739
$ws.Cells.Item(100, 16).Value = 560.0833
$ws.Cells.Item(100, 17).Value = "inf"
$ws.Cells.Item(100, 18).Value = "inf"
$ws.Cells.Item(100, 19).Value = "inf"
$ws.Cells.Item(100, 20).Value = "inf"
$ws.Cells.Item(100, 21).Value = "inf"
$ws.Cells.Item(100, 22).Value = "inf"
$ws.Cells.Item(100, 23).Value = "inf"
$ws.Cells.Item(100, 24).Value = "inf"
$ws.Cells.Item(100, 25).Value = "inf"
$ws.Cells.Item(100, 26).Value = "inf"
$ws.Cells.Item(100, 27).Value = "inf"
$ws.Cells.Item(100, 28).Value = "inf"
$ws.Cells.Item(100, 29).Value = "inf"
$ws.Cells.Item(100, 30).Value = "inf"
$ws.Cells.Item(100, 31).Value = "inf"
$ws.Cells.Item(100, 32).Value = "inf"
$ws.Cells.Item(100, 33).Value = "inf"
$ws.Cells.Item(100, 34).Value = "inf"
$ws.Cells.Item(100, 35).Value = "inf"
$ws.Cells.Item(100, 36).Value = "inf"
$ws.Cells.Item(100, 37).Value = "inf"
$ws.Cells.Item(100, 38).Value = "inf"
$ws.Cells.Item(100, 39).Value = "inf"
$ws.Cells.Item(100, 40).Value = "inf"
$ws.Cells.Item(100, 41).Value = "inf"
$ws.Cells.Item(100, 42).Value = "inf"
$ws.Cells.Item(100, 43).Value = "inf"
$ws.Cells.Item(100, 44).Value = "inf"
$ws.Cells.Item(100, 45).Value = "inf"

# Row 101
$ws.Cells.Item(101, 10).Value = 454694.771829416
$ws.Cells.Item(101, 11).Value = 488921.317581452
$ws.Cells.Item(101, 12).Value = 540670.793726652
$ws.Cells.Item(101, 13).Value = 580229.131052459
$ws.Cells.Item(101, 14).Value = 547339.903030655
$ws.Cells.Item(101, 15).Value = 545519.857413516
$ws.Cells.Item(101, 16).Value = 545519.9
$ws.Cells.Item(101, 17).Value = 573147.23119998
$ws.Cells.Item(101, 18).Value = 602173.722044243
$ws.Cells.Item(101, 19).Value = 632670.231628662
$ws.Cells.Item(101, 20).Value = 664711.2076399371
$ws.Cells.Item(101, 21).Value = 698374.868096333
$ws.Cells.Item(101, 22).Value = 733743.392292499
$ws.Cells.Item(101, 23).Value = 770903.121414501
$ws.Cells.Item(101, 24).Value = 809944.7693148199
$ws.Cells.Item(101, 25).Value = 850963.643961834
$ws.Cells.Item(101, 26).Value = 894059.880104412
$ws.Cells.Item(101, 27).Value = 939338.683719568
$ws.Cells.Item(101, 28).Value = 986910.588839938
$ws.Cells.Item(101, 29).Value = 1036891.72738804
$ws.Cells.Item(101, 30).Value = 1089404.11267604
$ws.Cells.Item(101, 31).Value = 1144575.93726305
$ws.Cells.Item(101, 32).Value = 1202541.88589717
$ws.Cells.Item(101, 33).Value = 1263443.46430619
$ws.Cells.Item(101, 34).Value = 1327429.3446395
$ws.Cells.Item(101, 35).Value = 1394655.72840465
$ws.Cells.Item(101, 36).Value = 1465286.72778448
$ws.Cells.Item(101, 37).Value = 1539494.76626564
$ws.Cells.Item(101, 38).Value = 1617460.99955659
$ws.Cells.Item(101, 39).Value = 1699375.7578226
$ws.Cells.Item(101, 40).Value = 1785439.0103173
$ws.Cells.Item(101, 41).Value = 1875860.85354502
$ws.Cells.Item(101, 42).Value = 1970862.0241457
$ws.Cells.Item(101, 43).Value = 2070674.43775433
$ws.Cells.Item(101, 44).Value = 2175541.75515039
$ws.Cells.Item(101, 45).Value = 2285719.97707947

# Row 102
$ws.Cells.Item(102, 10).Value = 46286.2291724209
$ws.Cells.Item(102, 11).Value = 52112.8840560455
$ws.Cells.Item(102, 12).Value = 73490.6204037136
$ws.Cells.Item(102, 13).Value = 78329.75859486101
$ws.Cells.Item(102, 14).Value = 54240.6933822808
$ws.Cells.Item(102, 15).Value = 46495.4777176006
$ws.Cells.Item(102, 16).Value = 46495.48
$ws.Cells.Item(102, 17).Value = 46391.7688979952
$ws.Cells.Item(102, 18).Value = 46288.2891301476
$ws.Cells.Item(102, 19).Value = 46185.0401804518
$ws.Cells.Item(102, 20).Value = 46082.0215340532
$ws.Cells.Item(102, 21).Value = 45979.232677246
$ws.Cells.Item(102, 22).Value = 45876.6730974698
$ws.Cells.Item(102, 23).Value = 45774.3422833078
$ws.Cells.Item(102, 24).Value = 45672.2397244839
$ws.Cells.Item(102, 25).Value = 45570.3649118601
$ws.Cells.Item(102, 26).Value = 45468.7173374342
$ws.Cells.Item(102, 27).Value = 45367.296494337
$ws.Cells.Item(102, 28).Value = 45266.1018768299
$ws.Cells.Item(102, 29).Value = 45165.1329803026
$ws.Cells.Item(102, 30).Value = 45064.3893012702
$ws.Cells.Item(102, 31).Value = 44963.8703373707
$ws.Cells.Item(102, 32).Value = 44863.575587363
$ws.Cells.Item(102, 33).Value = 44763.5045511238
$ws.Cells.Item(102, 34).Value = 44663.6567296454
$ws.Cells.Item(102, 35).Value = 44564.0316250331
$ws.Cells.Item(102, 36).Value = 44464.6287405031
$ws.Cells.Item(102, 37).Value = 44365.4475803793
$ws.Cells.Item(102, 38).Value = 44266.4876500916
$ws.Cells.Item(102, 39).Value = 44167.7484561727
$ws.Cells.Item(102, 40).Value = 44069.2295062562
$ws.Cells.Item(102, 41).Value = 43970.9303090742
$ws.Cells.Item(102, 42).Value = 43872.8503744541
$ws.Cells.Item(102, 43).Value = 43774.989213317
$ws.Cells.Item(102, 44).Value = 43677.3463376748
$ws.Cells.Item(102, 45).Value = 43579.9212606281

# Row 103
$ws.Cells.Item(103, 10).Value = 516717.294788736
$ws.Cells.Item(103, 11).Value = 509102.946985872
$ws.Cells.Item(103, 12).Value = 549636.164840267
$ws.Cells.Item(103, 13).Value = 536396.831873262
$ws.Cells.Item(103, 14).Value = 544549.298329128
$ws.Cells.Item(103, 15).Value = 552048.372609986
$ws.Cells.Item(103, 16).Value = 552048.4
$ws.Cells.Item(103, 17).Value = 570788.958169804
$ws.Cells.Item(103, 18).Value = 590165.707877371
$ws.Cells.Item(103, 19).Value = 610200.246114054
$ws.Cells.Item(103, 20).Value = 630914.903030964
$ws.Cells.Item(103, 21).Value = 652332.766827776
$ws.Cells.Item(103, 22).Value = 674477.709486437
$ws.Cells.Item(103, 23).Value = 697374.413378463
$ws.Cells.Item(103, 24).Value = 721048.39877549
$ws.Cells.Item(103, 25).Value = 745526.052293725
$ws.Cells.Item(103, 26).Value = 770834.656304016
$ws.Cells.Item(103, 27).Value = 797002.419340312
$ws.Cells.Item(103, 28).Value = 824058.5075404081
$ws.Cells.Item(103, 29).Value = 852033.077154019
$ws.Cells.Item(103, 30).Value = 880957.308154421
$ws.Cells.Item(103, 31).Value = 910863.438991104
$ws.Cells.Item(103, 32).Value = 941784.802522202
$ws.Cells.Item(103, 33).Value = 973755.863166712
$ws.Cells.Item(103, 34).Value = 1006812.25531795
$ws.Cells.Item(103, 35).Value = 1040990.82306103
$ws.Cells.Item(103, 36).Value = 1076329.66123864
$ws.Cells.Item(103, 37).Value = 1112868.1579109
$ws.Cells.Item(103, 38).Value = 1150647.03825663
$ws.Cells.Item(103, 39).Value = 1189708.40996491
$ws.Cells.Item(103, 40).Value = 1230095.81016759
$ws.Cells.Item(103, 41).Value = 1271854.25396505
$ws.Cells.Item(103, 42).Value = 1315030.28459923
$ws.Cells.Item(103, 43).Value = 1359672.0253299
$ws.Cells.Item(103, 44).Value = 1405829.23307211
$ws.Cells.Item(103, 45).Value = 1453553.35385427

# Row 104
$ws.Cells.Item(104, 10).Value = 192688.623563331
$ws.Cells.Item(104, 11).Value = 200550.794308594
$ws.Cells.Item(104, 12).Value = 207914.297062672
$ws.Cells.Item(104, 13).Value = 208783.213786659
$ws.Cells.Item(104, 14).Value = 218368.578380171
$ws.Cells.Item(104, 15).Value = 214829.978676142
$ws.Cells.Item(104, 16).Value = 214830
$ws.Cells.Item(104, 17).Value = 224258.757309972
$ws.Cells.Item(104, 18).Value = 234101.337011651
$ws.Cells.Item(104, 19).Value = 244375.901516714
$ws.Cells.Item(104, 20).Value = 255101.410373982
$ws.Cells.Item(104, 21).Value = 266297.655255273
$ws.Cells.Item(104, 22).Value = 277985.296476781
$ws.Cells.Item(104, 23).Value = 290185.90112334
$ws.Cells.Item(104, 24).Value = 302921.982845946
$ws.Cells.Item(104, 25).Value = 316217.04340597
$ws.Cells.Item(104, 26).Value = 330095.616042714
$ws.Cells.Item(104, 27).Value = 344583.31074435
$ws.Cells.Item(104, 28).Value = 359706.861505767
$ws.Cells.Item(104, 29).Value = 375494.175660538
$ws.Cells.Item(104, 30).Value = 391974.385378041
$ws.Cells.Item(104, 31).Value = 409177.90142075
$ws.Cells.Item(104, 32).Value = 427136.469260902
$ws.Cells.Item(104, 33).Value = 445883.227660098
$ws.Cells.Item(104, 34).Value = 465452.769819918
$ws.Cells.Item(104, 35).Value = 485881.2072164
$ws.Cells.Item(104, 36).Value = 507206.236236181
$ws.Cells.Item(104, 37).Value = 529467.207737253
$ws.Cells.Item(104, 38).Value = 552705.199662697
$ws.Cells.Item(104, 39).Value = 576963.092841393
$ws.Cells.Item(104, 40).Value = 602285.650115574
$ws.Cells.Item(104, 41).Value = 628719.598941243
$ws.Cells.Item(104, 42).Value = 656313.717613867
$ws.Cells.Item(104, 43).Value = 685118.92527847
$ws.Cells.Item(104, 44).Value = 715188.3758902079
$ws.Cells.Item(104, 45).Value = 746577.556298818

# Row 111
$ws.Cells.Item(111, 10).Value = 2894.36460075465
$ws.Cells.Item(111, 11).Value = 3189.68732307984
$ws.Cells.Item(111, 12).Value = 2345.54950698173
$ws.Cells.Item(111, 13).Value = 2153.27641877544
$ws.Cells.Item(111, 14).Value = 2334.17669621676
$ws.Cells.Item(111, 15).Value = 2175.66975644922
$ws.Cells.Item(111, 16).Value = 2175.67
$ws.Cells.Item(111, 17).Value = 2214.15404051064
$ws.Cells.Item(111, 18).Value = 2253.3188006957
$ws.Cells.Item(111, 19).Value = 2293.17632137179
$ws.Cells.Item(111, 20).Value = 2333.73885633791
$ws.Cells.Item(111, 21).Value = 2375.01887614264
$ws.Cells.Item(111, 22).Value = 2417.0290719182
$ws.Cells.Item(111, 23).Value = 2459.78235928214
$ws.Cells.Item(111, 24).Value = 2503.29188230814
$ws.Cells.Item(111, 25).Value = 2547.57101756703
$ws.Cells.Item(111, 26).Value = 2592.63337823928
$ws.Cells.Item(111, 27).Value = 2638.49281830023
$ws.Cells.Item(111, 28).Value = 2685.16343677937
$ws.Cells.Item(111, 29).Value = 2732.65958209493
$ws.Cells.Item(111, 30).Value = 2780.99585646519
$ws.Cells.Item(111, 31).Value = 2830.18712039774
$ws.Cells.Item(111, 32).Value = 2880.24849725824
$ws.Cells.Item(111, 33).Value = 2931.19537791993
$ws.Cells.Item(111, 34).Value = 2983.04342549538
$ws.Cells.Item(111, 35).Value = 3035.808580152
$ws.Cells.Item(111, 36).Value = 3089.50706401266
$ws.Cells.Item(111, 37).Value = 3144.15538614303
$ws.Cells.Item(111, 38).Value = 3199.77034762712
$ws.Cells.Item(111, 39).Value = 3256.36904673261
$ws.Cells.Item(111, 40).Value = 3313.96888416754
$ws.Cells.Item(111, 41).Value = 3372.58756843
$ws.Cells.Item(111, 42).Value = 3432.24312125242
$ws.Cells.Item(111, 43).Value = 3492.95388314216
$ws.Cells.Item(111, 44).Value = 3554.73851902015
$ws.Cells.Item(111, 45).Value = 3617.61602395918

# Row 112
$ws.Cells.Item(112, 10).Value = 135537.028549838
$ws.Cells.Item(112, 11).Value = 136312.736730274
$ws.Cells.Item(112, 12).Value = 138978.691114166
$ws.Cells.Item(112, 13).Value = 140670.504929068
$ws.Cells.Item(112, 14).Value = 134148.065059823
$ws.Cells.Item(112, 15).Value = 123085.108589292
$ws.Cells.Item(112, 16).Value = 123085.1
$ws.Cells.Item(112, 17).Value = 123498.408540439
$ws.Cells.Item(112, 18).Value = 123913.104933263
$ws.Cells.Item(112, 19).Value = 124329.193838754
$ws.Cells.Item(112, 20).Value = 124746.679932842
$ws.Cells.Item(112, 21).Value = 125165.567907159
$ws.Cells.Item(112, 22).Value = 125585.862469091
$ws.Cells.Item(112, 23).Value = 126007.568341832
$ws.Cells.Item(112, 24).Value = 126430.690264433
$ws.Cells.Item(112, 25).Value = 126855.232991862
$ws.Cells.Item(112, 26).Value = 127281.201295052
$ws.Cells.Item(112, 27).Value = 127708.599960955
$ws.Cells.Item(112, 28).Value = 128137.4337926
$ws.Cells.Item(112, 29).Value = 128567.707609141
$ws.Cells.Item(112, 30).Value = 128999.426245918
$ws.Cells.Item(112, 31).Value = 129432.594554504
$ws.Cells.Item(112, 32).Value = 129867.217402765
$ws.Cells.Item(112, 33).Value = 130303.299674914
$ws.Cells.Item(112, 34).Value = 130740.846271561
$ws.Cells.Item(112, 35).Value = 131179.862109776
$ws.Cells.Item(112, 36).Value = 131620.352123137
$ws.Cells.Item(112, 37).Value = 132062.32126179
$ws.Cells.Item(112, 38).Value = 132505.774492504
$ws.Cells.Item(112, 39).Value = 132950.716798723
$ws.Cells.Item(112, 40).Value = 133397.153180628
$ws.Cells.Item(112, 41).Value = 133845.088655188
$ws.Cells.Item(112, 42).Value = 134294.52825622
$ws.Cells.Item(112, 43).Value = 134745.477034443
$ws.Cells.Item(112, 44).Value = 135197.940057537
$ws.Cells.Item(112, 45).Value = 135651.922410198

